$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Group") to make room for "Age"
$ws.Columns.Item(2).Insert()

# Header for new Age column
$ws.Cells.Item(1, 2).Value = "Age"
$ws.Cells.Item(1, 2).Font.Bold = $true
$ws.Cells.Item(1, 2).HorizontalAlignment = $ws.Cells.Item(1, 1).HorizontalAlignment

# Age values for rows 2..26 (EHP01..EHP25)
$ages = @(47, 51, 72, 54, 55, 52, 63, 58, 56, 50, 45, 61, 73, 70, 63, 75, 70, 66, 65, 66, 68, 80, 66, 70, 67)

for ($i = 0; $i -lt $ages.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $ages[$i]
}
